$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename/add Programs/Homework/Literacy/Technology/Reading Specialist columns ---
$ws.Range("G1").Value = "Programs"
$ws.Range("H1").Value = "Homework"
$ws.Range("I1").Value = "Literacy"
$ws.Range("J1").Value = "Technology"
$ws.Range("K1").Value = "Reading Specialist"

# --- Column widths: col 7 shrinks to old col-8 width; cols 8-11 get new widths ---
# (input values are chosen so that Excel's internal character-width/pixel
# rounding lands on the intended stored width as closely as possible)
$ws.Columns.Item(7).ColumnWidth = 6.71525
$ws.Columns.Item(8).ColumnWidth = 9.00125
$ws.Columns.Item(9).ColumnWidth = 4.5725
$ws.Columns.Item(10).ColumnWidth = 11.144
$ws.Columns.Item(11).ColumnWidth = 12.28675

# --- Row 3: Active flag flips from TRUE to FALSE ---
$ws.Range("E3").Value = $false

# --- Data rows 2-54: column G becomes numeric "Programs" count (old H value),
#     columns H-K become new Homework/Literacy/Technology/Reading Specialist numbers ---
$data = @{
    2 = @(3,0,0,0,0)
    3 = @(0,0,0,0,0)
    4 = @(0,0,0,0,0)
    5 = @(0,0,0,0,0)
    6 = @(0,0,0,0,0)
    7 = @(0,0,0,0,0)
    8 = @(0,0,0,0,0)
    9 = @(3,15,20,25,30)
    10 = @(3,30,20,10,0)
    11 = @(3,40,10,20,15)
    12 = @(0,0,0,0,0)
    13 = @(0,0,0,0,0)
    14 = @(3,0,0,0,0)
    15 = @(0,0,0,0,0)
    16 = @(0,0,0,0,0)
    17 = @(0,0,0,0,0)
    18 = @(0,0,0,0,0)
    19 = @(3,15,8,5,20)
    20 = @(3,0,0,0,0)
    21 = @(3,0,0,0,0)
    22 = @(3,0,0,0,0)
    23 = @(3,0,0,0,0)
    24 = @(0,0,0,0,0)
    25 = @(0,0,0,0,0)
    26 = @(0,0,0,0,0)
    27 = @(0,0,0,0,0)
    28 = @(0,0,0,0,0)
    29 = @(0,0,0,0,0)
    30 = @(0,0,0,0,0)
    31 = @(0,0,0,0,0)
    32 = @(0,0,0,0,0)
    33 = @(0,0,0,0,0)
    34 = @(0,0,0,0,0)
    35 = @(0,0,0,0,0)
    36 = @(0,0,0,0,0)
    37 = @(0,0,0,0,0)
    38 = @(0,0,0,0,0)
    39 = @(0,0,0,0,0)
    40 = @(0,0,0,0,0)
    41 = @(0,0,0,0,0)
    42 = @(0,0,0,0,0)
    43 = @(0,0,0,0,0)
    44 = @(0,0,0,0,0)
    45 = @(0,0,0,0,0)
    46 = @(0,0,0,0,0)
    47 = @(0,0,0,0,0)
    48 = @(0,0,0,0,0)
    49 = @(0,0,0,0,0)
    50 = @(0,0,0,0,0)
    51 = @(0,0,0,0,0)
    52 = @(0,0,0,0,0)
    53 = @(0,0,0,0,0)
    54 = @(0,0,0,0,0)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 7).Value = $vals[0]
    $ws.Cells.Item($r, 8).Value = $vals[1]
    $ws.Cells.Item($r, 9).Value = $vals[2]
    $ws.Cells.Item($r, 10).Value = $vals[3]
    $ws.Cells.Item($r, 11).Value = $vals[4]
}
